# Apply the "finish filtering of all benchmarks" edit:
#  - update a handful of rank (G) values and two note (K) texts
#  - highlight several "title" (C) cells with the built-in Good / Neutral
#    cell styles to flag filtered / reviewed rows
#  - move the active selection from G11 to F10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rank (column G) corrections ---------------------------------------
$ws.Range("G4").Value  = 3
$ws.Range("G5").Value  = 3
$ws.Range("G8").Value  = 4
$ws.Range("G13").Value = 2
$ws.Range("G17").Value = 2
$ws.Range("G21").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("G24").Value = 0

# --- note (column K) text updates ---------------------------------------
$ws.Range("K4").Value  = "stapled peptide; FITC group at C terminal for fluorescently labeled"
$ws.Range("K11").Value = "constrained peptide; FITC group for fluorescently labeled"

# --- highlight reviewed rows with built-in cell styles -------------------
$ws.Range("C2").Style  = "Good"
$ws.Range("C3").Style  = "Good"
$ws.Range("C6").Style  = "Neutral"
$ws.Range("C7").Style  = "Good"
$ws.Range("C8").Style  = "Good"
$ws.Range("C12").Style = "Good"
$ws.Range("C14").Style = "Good"
$ws.Range("C15").Style = "Good"
$ws.Range("C16").Style = "Good"

# --- move the active selection -------------------------------------------
$ws.Range("F10").Select()
